# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update MyForecast (column D) values ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("D3").Value = 52
$wsForecast.Range("D5").Value = 53
$wsForecast.Range("D7").Value = 55
$wsForecast.Range("D8").Value = 53
$wsForecast.Range("D10").Value = 51
$wsForecast.Range("D11").Value = 46
$wsForecast.Range("D12").Value = 41
$wsForecast.Range("D13").Value = 34

# --- Sheet "Summary": update forecast totals (column B) values ---
# These are stored as text (not numbers), so force a text number format
# before assigning the numeric-looking strings, otherwise Excel would
# auto-convert them back to numeric cells.
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9:B11").NumberFormat = "@"
$wsSummary.Range("B9").Value = "797"
$wsSummary.Range("B10").Value = "428"
$wsSummary.Range("B11").Value = "216"
